$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Login with valid username and password"
$ws.Range("B4").Value = "PASSED"
$ws.Range("C4").Value = "chrome"
